$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (header/selection row): SchoolName / Classroom Name / Section Name
# These are plain alphanumeric text values, assigning directly keeps them as text.
$ws.Cells.Item(2,1).Value = "FPK12School82076"
$ws.Cells.Item(2,2).Value = "FPK12Classroom82865"
$ws.Cells.Item(2,3).Value = "FPK12Section39584"

# E3/E4/E5 hold purely-numeric-looking text (e.g. "67807") that must stay text,
# not be auto-converted to a Number by the smart-typing in Range.Value. Stage the
# value in a scratch cell formatted as Text, copy it, and paste-special (values
# only) into the destination so the destination keeps its own existing style.
$scratch = $ws.Cells.Item(200, 200)

$scratch.NumberFormat = "@"
$scratch.Value = "67807"
$scratch.Copy()
$ws.Cells.Item(3,5).PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "21227"
$scratch.Copy()
$ws.Cells.Item(4,5).PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "34255"
$scratch.Copy()
$ws.Cells.Item(5,5).PasteSpecial(-4163)

$scratch.Clear()

$excel.CutCopyMode = 0
